$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("W1").Value = 0.81535464425547022
$ws.Range("AJ1").Value = 0.97347934491563315
$ws.Range("A2").Value = 0.79711416233203547
$ws.Range("J2").Value = 0.91921606069174344
$ws.Range("AS2").Value = 0.78548014363600926
$ws.Range("BC3").Value = 0.97411166095329427
$ws.Range("AT4").Value = 0.81875790657938696
$ws.Range("BN4").Value = 0.53861966207184131
$ws.Range("F5").Value = 0.77457588789477128
$ws.Range("X6").Value = 0.92157000195238836
$ws.Range("M7").Value = 0.8070067064894666
$ws.Range("J8").Value = 0.91236848701384599
$ws.Range("AB8").Value = 0.63610433852136228
$ws.Range("AO8").Value = 0.66659347267301694
$ws.Range("AY8").Value = 0.93609005650059385
$ws.Range("M11").Value = 0.80783009190845556
$ws.Range("AA11").Value = 0.83550837302522485
$ws.Range("G12").Value = 0.95362879164574754
$ws.Range("I12").Value = 0.80115907109769968
$ws.Range("N13").Value = 0.99567343582905521
$ws.Range("AC13").Value = 0.60732205231760528
$ws.Range("Y14").Value = 0.95876369124220195
$ws.Range("AI14").Value = 0.91209878751200213
$ws.Range("BP14").Value = 0.96095529229039101
$ws.Range("AQ15").Value = 0.77007724588965354
$ws.Range("F16").Value = 0.91136965025582783
$ws.Range("AD16").Value = 0.89084579311904133
$ws.Range("S17").Value = 0.64875165885480912
$ws.Range("G18").Value = 0.99400871816479164
$ws.Range("AL18").Value = 0.9173344802154495
$ws.Range("BF18").Value = 0.968430666934049
$ws.Range("BJ19").Value = 0.85208756597441271
$ws.Range("Q20").Value = 0.97089824070128317
$ws.Range("I21").Value = 0.7887321342714122
$ws.Range("W21").Value = 0.79088610171581952
$ws.Range("AI21").Value = 0.84704485433920063
$ws.Range("AZ21").Value = 0.98158201391922639
$ws.Range("BA22").Value = 0.66972036291955117
$ws.Range("BH22").Value = 0.89951687238219669
$ws.Range("P23").Value = 0.75524975803487449
$ws.Range("AK23").Value = 0.98513809722119183
$ws.Range("BF24").Value = 0.73677271243919695
$ws.Range("BF25").Value = 0.69562366665212283
$ws.Range("C26").Value = 0.97709917399788437
$ws.Range("BM27").Value = 0.88155466523869874
$ws.Range("BP27").Value = 0.9308514287297085
$ws.Range("B28").Value = 0.77524672153854124
$ws.Range("E29").Value = 0.62404132637093923
$ws.Range("AF29").Value = 0.98114959586775363
$ws.Range("AJ29").Value = 0.93331296603816161
$ws.Range("W30").Value = 0.82367796973524232
$ws.Range("Z30").Value = 0.97983802265535491
$ws.Range("AT30").Value = 0.7989475780678017
$ws.Range("AF31").Value = 0.86303267917034066
$ws.Range("AK31").Value = 0.81714969314354113
$ws.Range("X32").Value = 0.93941499679131479
$ws.Range("AD32").Value = 0.91186268242651458
$ws.Range("BF32").Value = 0.6654991762088921
$ws.Range("AC34").Value = 0.9626305082890978
$ws.Range("AR34").Value = 0.96051604029942472
$ws.Range("AJ35").Value = 0.66931847650178256
$ws.Range("D36").Value = 0.96628045553914821
$ws.Range("L36").Value = 0.90606253857766639
$ws.Range("U37").Value = 0.99196140191755333
$ws.Range("AT37").Value = 0.63941219450110975
$ws.Range("BN37").Value = 0.92307064435511199
$ws.Range("AG39").Value = 0.81006783074807243
$ws.Range("AV39").Value = 0.89474089559511394
$ws.Range("AL40").Value = 0.86473656645212049
$ws.Range("AY40").Value = 0.97699711584378335
$ws.Range("AK41").Value = 0.90986066338413507
$ws.Range("AL41").Value = 0.97804628128114324
$ws.Range("AB43").Value = 0.87807296774601218
$ws.Range("AR43").Value = 0.82472624672699713
$ws.Range("B44").Value = 0.95859703601678059
$ws.Range("AP44").Value = 0.98644139570324574
$ws.Range("E45").Value = 0.71883074053721607
$ws.Range("H45").Value = 0.76777978395652047
$ws.Range("Y45").Value = 0.8544689743619226
$ws.Range("AG45").Value = 0.9895998824031117
$ws.Range("AR45").Value = 0.86413025922305953
$ws.Range("BB45").Value = 0.9252603756643385
$ws.Range("AU46").Value = 0.84323132305124948
$ws.Range("B47").Value = 0.89972697090764631
$ws.Range("BF47").Value = 0.65684526638992979
$ws.Range("I48").Value = 0.91801708938050197
$ws.Range("S48").Value = 0.54297519733390298
$ws.Range("Y48").Value = 0.7920613922552413
$ws.Range("AK48").Value = 0.93315715699739377
$ws.Range("AE49").Value = 0.88036959301212026
$ws.Range("AU49").Value = 0.91331318490406965
$ws.Range("BD49").Value = 0.81739025293621448
$ws.Range("BH50").Value = 0.7356301385179167
$ws.Range("O51").Value = 0.99328437247064194
$ws.Range("BE51").Value = 0.96446061576584519
$ws.Range("AX52").Value = 0.59213036670090169
$ws.Range("BF52").Value = 0.82691054361558192
$ws.Range("C53").Value = 0.846149787624494
$ws.Range("AN53").Value = 0.77662255788160994
$ws.Range("BB53").Value = 0.71143369429290026
$ws.Range("BH53").Value = 0.82024825285021796
$ws.Range("G54").Value = 0.75369156733482823
$ws.Range("BC54").Value = 0.78651114755187923
$ws.Range("T55").Value = 0.73230425592780302
$ws.Range("AI55").Value = 0.8932630508161381
$ws.Range("AH56").Value = 0.93771025485342963
$ws.Range("AY56").Value = 0.99843079138180424
$ws.Range("AP57").Value = 0.81998322831438464
$ws.Range("O58").Value = 0.57065664131384031
$ws.Range("S59").Value = 0.99908920495838149
$ws.Range("U59").Value = 0.68496088495999852
$ws.Range("AM59").Value = 0.86867721846597501
$ws.Range("BN59").Value = 0.81388246485614646
$ws.Range("BP59").Value = 0.76381044835924561
$ws.Range("L60").Value = 0.87442207690629026
$ws.Range("M60").Value = 0.9702420309743498
$ws.Range("R60").Value = 0.93715534524558231
$ws.Range("K61").Value = 0.65001047291403102
$ws.Range("AI61").Value = 0.80272736476224005
$ws.Range("BL62").Value = 0.87750554399675484
$ws.Range("AJ63").Value = 0.84785549353629464
$ws.Range("AT63").Value = 0.82601359657063322
$ws.Range("AY63").Value = 0.72386580975501036
$ws.Range("B64").Value = 0.64880946237965209
$ws.Range("F64").Value = 0.8109421555121441
$ws.Range("Q64").Value = 0.81058577316638036
$ws.Range("BN64").Value = 0.80146774345419369
$ws.Range("AT65").Value = 0.81443164124754208
$ws.Range("BO65").Value = 0.95973124279725175
$ws.Range("Z66").Value = 0.93089599663645528
$ws.Range("J67").Value = 0.93921973344481791
$ws.Range("L67").Value = 0.92863781701922543
$ws.Range("X67").Value = 0.95427193877013117
$ws.Range("L68").Value = 0.6672054086020055
$ws.Range("BE68").Value = 0.99559948718836644
$ws.Range("BL68").Value = 0.90042294868225792

# Column 62 (BJ) width shrinks from 12.7109375 to 11.7109375 character units.
# ColumnWidth snaps to the engine's internal pixel grid (steps of 1/6 char), so
# 10.83 is the input that lands on the closest achievable width (~11.6667).
$ws.Columns("BJ:BJ").ColumnWidth = 10.83
